$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Edit row 3 (was the "SIPL5316/SIPL5688/Doc Retrieval" order, formerly row 6) ---
$ws.Range("A3").Value = 45440.041666608799
$ws.Range("C3").Value = "SIPL5316"
$ws.Range("D3").Value = "SIPL5688"
$ws.Range("E3").Clear()
$ws.Range("F3").Clear()
$ws.Range("I3").Value = "Search"
$ws.Range("J3").Value = "Doc Retrieval"
$ws.Range("L3").Value = "Autauga"
$ws.Range("N3").Value = "Search(T1)"

# --- Edit row 4 (the "Typing" order, date now 45443) ---
$ws.Range("A4").Value = 45443.041666608799
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("I4").Value = "Typing"
$ws.Range("J4").Value = "Typing"
$ws.Range("L4").Value = "Autauga"
$ws.Range("M4").Value = "Typing"
$ws.Range("N4").Value = "Typing(T1)"

# --- Remove the now-stale trailing rows 5-9 ---
$ws.Rows("5:9").Delete()

# Cosmetic: update the selection to match the author's final cursor position
[void]$ws.Range("I8").Select()
